# Apply updated coin price/volume data as described in the commit diff.
# Numeric-looking text values (Price/Volume columns) are written with a
# leading apostrophe so Excel keeps them as literal text (matching the
# original inlineStr cell contents) instead of auto-converting them to
# numbers or percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''301.71'
$ws.Range("E2").Value = '''-1.03%'
$ws.Range("D3").Value = '''31.38'
$ws.Range("E3").Value = '''-2.27%'
$ws.Range("D4").Value = '''5.124'
$ws.Range("E4").Value = '''-2.75%'
$ws.Range("D5").Value = '''0.07385'
$ws.Range("E5").Value = '''-2.40%'
$ws.Range("D6").Value = '''2.130'
$ws.Range("E6").Value = '''33.81%'
$ws.Range("D7").Value = '''7.918'
$ws.Range("E7").Value = '''0.94%'
$ws.Range("D8").Value = '''3.827'
$ws.Range("E8").Value = '''-0.82%'
$ws.Range("D9").Value = '''0.9195'
$ws.Range("E9").Value = '''-0.93%'
$ws.Range("D10").Value = '''0.1706'
$ws.Range("E10").Value = '''1.29%'
$ws.Range("D11").Value = '''0.07483'
$ws.Range("E11").Value = '''-6.74%'
$ws.Range("D12").Value = '''0.08154'
$ws.Range("E12").Value = '''1.76%'
$ws.Range("D13").Value = '''0.03036'
$ws.Range("E13").Value = '''-0.15%'
$ws.Range("D14").Value = '''0.09917'
$ws.Range("E14").Value = '''-0.19%'
$ws.Range("E15").Value = '''-0.03%'
$ws.Range("D16").Value = '''0.006086'
$ws.Range("E16").Value = '''-2.90%'
$ws.Range("D17").Value = '''3.473'
$ws.Range("E17").Value = '''0.75%'
$ws.Range("E18").Value = '''-0.53%'
$ws.Range("D19").Value = '''0.3267'
$ws.Range("E19").Value = '''-1.01%'
$ws.Range("E20").Value = '''-1.18%'
$ws.Range("D21").Value = '''4.645'
$ws.Range("E21").Value = '''2.22%'
$ws.Range("D22").Value = '''0.04644'
$ws.Range("E22").Value = '''1.15%'
$ws.Range("D23").Value = '''0.1567'
$ws.Range("E23").Value = '''-3.15%'
$ws.Range("D24").Value = '''0.001225'
$ws.Range("E24").Value = '''1.08%'
$ws.Range("D25").Value = '''0.004484'
$ws.Range("E25").Value = '''0.10%'
$ws.Range("D26").Value = '''0.0001298'
$ws.Range("E26").Value = '''-6.89%'
$ws.Range("D27").Value = '''0.0003425'
$ws.Range("E27").Value = '''92.25%'
$ws.Range("D39").Value = '''0.01736'
$ws.Range("E39").Value = '''0.78%'
$ws.Range("D40").Value = '''0.04501'
$ws.Range("E40").Value = '''-0.16%'
$ws.Range("D41").Value = '''0.007368'
$ws.Range("E41").Value = '''6.69%'
$ws.Range("D42").Value = '''0.1348'
$ws.Range("E42").Value = '''-0.72%'
$ws.Range("D43").Value = '''0.002226'
$ws.Range("E43").Value = '''7.50%'
$ws.Range("E44").Value = '''-23.01%'
$ws.Range("D45").Value = '''0.00006291'
$ws.Range("E45").Value = '''2.22%'
$ws.Range("B46").Value = 'CoinbaseStockToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D46").Value = '''0.009992'
$ws.Range("E46").Value = '''-22.98%'
$ws.Range("B47").Value = 'BOLO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D47").Value = '''0.8085'
$ws.Range("E47").Value = '''14.02%'
